# Apply scheduled runner updates to Ravana_Profits leve-profit workbook
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("WVR")

# ALC sheet value updates
$ws1.Range("H18").Value = 273.7
$ws1.Range("I18").Value = 273.7
$ws1.Range("K18").Value = 273.7
$ws1.Range("M18").Value = 10.30000000000001
$ws1.Range("H29").Value = 3989.6667
$ws1.Range("J29").Value = 3989.6667
$ws1.Range("L29").Value = 11969.0001
$ws1.Range("N29").Value = -12531.0001
$ws1.Range("H38").Value = 2454.6667
$ws1.Range("I38").Value = 427.5
$ws1.Range("J38").Value = 6509
$ws1.Range("K38").Value = 1282.5
$ws1.Range("L38").Value = 19527
$ws1.Range("M38").Value = -910.5
$ws1.Range("N38").Value = -20271
$ws1.Range("H43").Value = 6942.3335
$ws1.Range("I43").Value = 6663.5
$ws1.Range("K43").Value = 6663.5
$ws1.Range("M43").Value = -6594.5
$ws1.Range("H53").Value = 321
$ws1.Range("I53").Value = 317.2
$ws1.Range("K53").Value = 317.2
$ws1.Range("M53").Value = 319.8
$ws1.Range("H58").Value = 7286.1665
$ws1.Range("J58").Value = 10729.25
$ws1.Range("L58").Value = 32187.75
$ws1.Range("N58").Value = -32487.75
$ws1.Range("H86").Value = 1428.3334
$ws1.Range("I86").Value = 1428.3334
$ws1.Range("K86").Value = 1428.3334
$ws1.Range("M86").Value = -305.3334
$ws1.Range("H89").Value = 1428.3334
$ws1.Range("I89").Value = 1428.3334
$ws1.Range("K89").Value = 7141.666999999999
$ws1.Range("M89").Value = -1525.666999999999
$ws1.Range("H101").Value = 1066.6666
$ws1.Range("I101").Value = 1066.6666
$ws1.Range("K101").Value = 3199.9998
$ws1.Range("M101").Value = -1577.9998
$ws1.Range("H116").Value = 9165.091
$ws1.Range("I116").Value = 8822.166999999999
$ws1.Range("J116").Value = 9576.6
$ws1.Range("K116").Value = 8822.166999999999
$ws1.Range("L116").Value = 9576.6
$ws1.Range("M116").Value = -5380.166999999999
$ws1.Range("N116").Value = -16460.6
$ws1.Range("H125").Value = 604.75
$ws1.Range("I125").Value = 134
$ws1.Range("J125").Value = 761.6667
$ws1.Range("K125").Value = 1206
$ws1.Range("L125").Value = 6855.0003
$ws1.Range("M125").Value = 1254
$ws1.Range("N125").Value = -11775.0003
$ws1.Range("H132").Value = 1207.28
$ws1.Range("I132").Value = 1207.28
$ws1.Range("K132").Value = 3621.84
$ws1.Range("M132").Value = -1091.84
$ws1.Range("H137").Value = 2059.4783
$ws1.Range("I137").Value = 1578.8462
$ws1.Range("K137").Value = 4736.5386
$ws1.Range("M137").Value = -2186.5386

# ARM sheet value updates
$ws2.Range("H2").Value = 1095.9048
$ws2.Range("I2").Value = 1025.9286
$ws2.Range("J2").Value = 1235.8572
$ws2.Range("K2").Value = 1025.9286
$ws2.Range("L2").Value = 1235.8572
$ws2.Range("M2").Value = -912.9286
$ws2.Range("N2").Value = -1461.8572
$ws2.Range("H102").Value = 2446.2
$ws2.Range("I102").Value = 2218
$ws2.Range("K102").Value = 2218
$ws2.Range("M102").Value = -596
$ws2.Range("H116").Value = 1095.9048
$ws2.Range("I116").Value = 1025.9286
$ws2.Range("J116").Value = 1235.8572
$ws2.Range("K116").Value = 1025.9286
$ws2.Range("L116").Value = 1235.8572
$ws2.Range("M116").Value = 1268.0714
$ws2.Range("N116").Value = -5823.8572
$ws2.Range("H122").Value = 5554
$ws2.Range("I122").Value = 5776
$ws2.Range("K122").Value = 17328
$ws2.Range("M122").Value = -14878

# BSM sheet value updates
$ws3.Range("H3").Value = 1095.9048
$ws3.Range("I3").Value = 1025.9286
$ws3.Range("J3").Value = 1235.8572
$ws3.Range("K3").Value = 1025.9286
$ws3.Range("L3").Value = 1235.8572
$ws3.Range("M3").Value = -911.9286
$ws3.Range("N3").Value = -1463.8572
$ws3.Range("H22").Value = 715
$ws3.Range("I22").Value = 715
$ws3.Range("K22").Value = 715
$ws3.Range("M22").Value = -542
$ws3.Range("H134").Value = 2485.2144
$ws3.Range("I134").Value = 2174.875
$ws3.Range("K134").Value = 6524.625
$ws3.Range("M134").Value = -3989.625

# CRP sheet value updates
$ws4.Range("H7").Value = 288
$ws4.Range("I7").Value = 288
$ws4.Range("K7").Value = 288
$ws4.Range("M7").Value = -175
$ws4.Range("H31").Value = 1207.0625
$ws4.Range("I31").Value = 1118.75
$ws4.Range("K31").Value = 1118.75
$ws4.Range("M31").Value = -823.75
$ws4.Range("H34").Value = 1207.0625
$ws4.Range("I34").Value = 1118.75
$ws4.Range("K34").Value = 1118.75
$ws4.Range("M34").Value = -916.75
$ws4.Range("H134").Value = 3145
$ws4.Range("I134").Value = 3145
$ws4.Range("K134").Value = 9435
$ws4.Range("M134").Value = -6900

# CUL sheet value updates
$ws5.Range("H86").Value = 700
$ws5.Range("I86").Value = 700
$ws5.Range("K86").Value = 2100
$ws5.Range("M86").Value = -914
$ws5.Range("H89").Value = 700
$ws5.Range("I89").Value = 700
$ws5.Range("K89").Value = 6300
$ws5.Range("M89").Value = -372

# GSM sheet value updates
$ws6.Range("H132").Value = 2979.4375
$ws6.Range("I132").Value = 2084.75
$ws6.Range("K132").Value = 6254.25
$ws6.Range("M132").Value = -3724.25

# WVR sheet value updates
$ws7.Range("H132").Value = 4592.636
$ws7.Range("I132").Value = 4131.75
$ws7.Range("J132").Value = 4856
$ws7.Range("K132").Value = 12395.25
$ws7.Range("L132").Value = 14568
$ws7.Range("M132").Value = -9865.25
$ws7.Range("N132").Value = -19628
$ws7.Range("H136").Value = 2898
$ws7.Range("J136").Value = 3289.8
$ws7.Range("L136").Value = 9869.400000000001
$ws7.Range("N136").Value = -14969.4

# Row-level updates involving added/removed cells
# BSM row 76
$ws3.Range("H76").Value = 29999
$ws3.Range("I76").Value = 0
$ws3.Range("J76").Value = 29999
$ws3.Range("K76").Value = 0
$ws3.Range("L76").Value = 29999
$ws3.Range("N76").Value = -30629
$ws3.Range("M76").ClearContents()

# BSM row 79
$ws3.Range("H79").Value = 29999
$ws3.Range("I79").Value = 0
$ws3.Range("J79").Value = 29999
$ws3.Range("K79").Value = 0
$ws3.Range("L79").Value = 29999
$ws3.Range("N79").Value = -32183
$ws3.Range("M79").ClearContents()

# BSM row 99
$ws3.Range("H99").Value = 1019.6
$ws3.Range("J99").Value = 1100
$ws3.Range("L99").Value = 1100
$ws3.Range("N99").Value = -4096

# BSM row 102
$ws3.Range("H102").Value = 5000
$ws3.Range("I102").Value = 5000
$ws3.Range("J102").Value = 0
$ws3.Range("K102").Value = 5000
$ws3.Range("L102").Value = 0
$ws3.Range("M102").Value = -1755
$ws3.Range("N102").ClearContents()

# CRP row 105
$ws4.Range("H105").Value = 2498.5
$ws4.Range("J105").Value = 0
$ws4.Range("L105").Value = 0
$ws4.Range("N105").ClearContents()

# CUL row 21
$ws5.Range("H21").Value = 100
$ws5.Range("J21").Value = 100
$ws5.Range("L21").Value = 300
$ws5.Range("N21").Value = -646

# GSM row 107
$ws6.Range("H107").Value = 500
$ws6.Range("I107").Value = 500
$ws6.Range("K107").Value = 500
$ws6.Range("M107").Value = 1420

# GSM row 126
$ws6.Range("H126").Value = 1500
$ws6.Range("I126").Value = 1500
$ws6.Range("K126").Value = 4500
$ws6.Range("M126").Value = -2030
